$d = $word.ActiveDocument

function Replace-SubtitleText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found -and $rng.Find.Found) {
        $rng.Text = $newText
    }
}

Replace-SubtitleText 'known as the ants puzzle, which I''m' 'connus sous le nom « problème des fourmis », que je vais'
Replace-SubtitleText 'probably going to discuss in a different' 'probablement discuter dans une autre'
Replace-SubtitleText 'video. Let me just finish writing down' 'vidéo. Permettez-moi de finir d''écrire'
Replace-SubtitleText 'the title and, well, I can even draw a' 'le titre et, bien, je peux même dessiner une'
Replace-SubtitleText 'little ant right here. okay, let''s get' 'petite fourmi ici même. Bon, allons-y !'
Replace-SubtitleText 'started! As I said I''m going to discuss' 'allons-y ! Comme je disais avant, je vais discuter'
Replace-SubtitleText 'two puzzles in the first puzzle there' 'deux énigmes. Dans le premier,'
Replace-SubtitleText 'are two ants on a very high stool: a sort' 'il y a deux fourmis sur un plateau très haut, un genre'
Replace-SubtitleText 'of Mountain, flat on the top with two' 'de montagne, plat au sommet avec deux'
Replace-SubtitleText 'steep cliffs to both the sides. The flat' 'falaises abruptes sur les deux côtés. Le sommet'
Replace-SubtitleText 'peak is one meter wide the two ants move' 'plat mesure 1 mètre de large. Les deux fourmis bougent'
Replace-SubtitleText 'with a velocity, let''s call it V, which is' 'avec une vélocité, appelons-la v, qui est'
Replace-SubtitleText 'the same for both of them and that is' 'la même pour les deux et égal à'
Replace-SubtitleText 'equal to one centimeter per second. You' 'un centimètre par seconde. Vous'
Replace-SubtitleText 'can decide the direction towards each' 'pouvez décider la direction dans laquelle chaqu''une des'
Replace-SubtitleText 'ant moves if it is right or left and' 'la fourmi se déplace, soit vers la droite, soit vers la gauche,'
Replace-SubtitleText 'where exactly to place the two ants on the' 'et où placer exactement les deux fourmis'
Replace-SubtitleText 'top of the mountain. Your purpose is to' 'en haut de la montagne. Votre but est de'
